$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 444, shifting existing rows 444:478 down to 445:479.
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row 444 with the new daily price record.
$ws.Range("A444").Value = 8
$ws.Range("B444").Value = "Terminal La Palmera de La Serena"
$ws.Range("C444").Value = "Coquimbo"
$ws.Range("D444").Value = 45013
$ws.Range("E444").Value = 4
$ws.Range("F444").Value = 100112032
$ws.Range("G444").Value = "Zapallo italiano"
$ws.Range("H444").Value = "Sin especificar"
$ws.Range("I444").Value = "Primera"
$ws.Range("J444").Value = 300
$ws.Range("K444").Value = 9000
$ws.Range("L444").Value = 10000
$ws.Range("M444").Value = 9500
$ws.Range("N444").Value = "$/caja 70 unidades"
$ws.Range("O444").Value = "Provincia de Limarí"
$ws.Range("P444").Value = 136
$ws.Range("Q444").Value = 70
$ws.Range("R444").Value = "Hortaliza"
